# The commit swaps the contents of ppt/theme/theme1.xml (was "Office Theme")
# and ppt/theme/theme2.xml (was "Integral" / "Red Violet"): after the edit,
# theme1.xml holds the Integral/Red Violet theme and theme2.xml holds the
# plain Office Theme. Both files already share an identical <a:fontScheme>
# and <a:fmtScheme>; the two themes only differ in their <a:clrScheme>
# (the 12 theme colours) and the cosmetic name= attributes.
#
# ppt/theme/theme2.xml is the theme that the presentation's slide master
# (and therefore the whole deck) actually uses, so we recolor it in place
# to the Office Theme palette -- this reproduces the half of the swap that
# is visible to the document's live/active theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme

# Rename theme + colour scheme to match the incoming "Office Theme" content.
$theme.Name = "Office Theme"
$master.ColorScheme.Name = "Office"
$p.Designs.Item(1).Name = "Office Theme"

$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0x000000   # dk1
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1
$colors.Colors(3).RGB  = 0x6A5444   # dk2      (#44546A)
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      (#E7E6E6)
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  (#5B9BD5)
$colors.Colors(6).RGB  = 0x317DED   # accent2  (#ED7D31)
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  (#A5A5A5)
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  (#FFC000)
$colors.Colors(9).RGB  = 0xC47244   # accent5  (#4472C4)
$colors.Colors(10).RGB = 0x47AD70   # accent6  (#70AD47)
$colors.Colors(11).RGB = 0xC16305   # hlink    (#0563C1)
$colors.Colors(12).RGB = 0x724F95   # folHlink (#954F72)
